$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new training-run row (row 7) with the same shape as the
# existing data rows, continuing the log of "Bag" method runs.
$ws.Range("A7").Value = 42604.891446759262
$ws.Range("B7").Value = "Bag"
$ws.Range("C7").Value = 42
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
